# Add I0 and IF columns (I and J) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, mirroring the style of the existing header row (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Data rows 2-27: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 27; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2

    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
